# Auto-generated Excel COM-interop script
# Applies per-cell value updates to the crypto price/volume tracking sheet
# (column D = Price, column E = Volume(1h), column B/C = Coin name/Link)
# as captured by the authoritative OOXML diff for this commit.
#
# All changed cells in this sheet are stored as plain text (inlineStr) in
# the workbook -- e.g. "69.642.18" or "0.0452" -- not as numbers. Excel's
# COM layer will happily "helpfully" reinterpret a numeric-looking string
# assigned to .Value as an actual floating point number (losing trailing
# zeros / exact formatting / thousand-separator dots in the process), so
# we explicitly force each target cell to the Text number format ("@")
# before writing the new value. This guarantees the value round-trips as
# the exact text the source data expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.642.18'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.505.49'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.05'
$ws.Range('E5').Value = '  -1.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '194.82'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -5.60%  '
$ws.Range('E10').Value = '  +0.39%  '
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('E12').Value = '  -2.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.49'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.063.35'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '593.56'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.801.85'
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.07'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.123'
$ws.Range('E19').Value = '  +2.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.502.15'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.987'
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.20'
$ws.Range('E22').Value = '  +6.29%  '
$ws.Range('E23').Value = '  +2.80%  '
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '101.45'
$ws.Range('E25').Value = '  -3.97%  '
$ws.Range('E26').Value = '  +3.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.83'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.50'
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.16'
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('E30').Value = '  +1.36%  '
$ws.Range('E31').Value = '  +2.85%  '
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.10'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0822'
$ws.Range('E35').Value = '  +6.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.723.86'
$ws.Range('E36').Value = '  +2.84%  '
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('E40').Value = '  -0.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.28'
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '485.29'
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('E43').Value = '  -2.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0452'
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.81'
$ws.Range('E46').Value = '  -3.45%  '
$ws.Range('E47').Value = '  -1.14%  '
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('E49').Value = '  -4.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000246'
$ws.Range('E50').Value = '  +2.65%  '
$ws.Range('E51').Value = '  +10.30%  '
